$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 5 ("Anticipated Changes") - Content Placeholder gets new bullets
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(1)
$tr5 = $sh5.TextFrame.TextRange

$tr5.Text = "Additional design features for debugger`rMethods for setting break points, displaying variables, etc.`rAdd details for connection`rEverything the nxt needs to create a connection`rAdd details for timer`rResearch java timer classes/methods`rPossibly add fields to Message Handler Class"

$tr5.Paragraphs(2,1).IndentLevel = 2
$tr5.Paragraphs(4,1).IndentLevel = 2
$tr5.Paragraphs(6,1).IndentLevel = 2

$para5_4 = $tr5.Paragraphs(4,1)
$para5_4.Characters(1,15).Text = "Everything the "
$para5_4.Characters(16,3).Text = "nxt"
$para5_4.Characters(19,31).Text = " needs to create a connection"

# ---------------------------------------------------------------------
# Slide 6 ("Management Information") - Content Placeholder text + move
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(1)

$sh6.Left = 36.0
$sh6.Top = 116.64000701904297
$sh6.Width = 648.0
$sh6.Height = 363.3600158691406

$tr6 = $sh6.TextFrame.TextRange
$tr6.Text = "Time spent on design`rApprox. 2 hours discussing high level design`rProblems`rSending multiple messages/acknowledgements`rDecoding messages`rMajor Risks`rShould we allow sending multiple messages before receiving an acknowledgment?`rHow to decode message in Message Handler Class for use in Driver class to implement action`rDesign details of debugger " + [char]0x2013 + " how to set breakpoints "

$tr6.Paragraphs(2,1).IndentLevel = 2
$tr6.Paragraphs(4,1).IndentLevel = 2
$tr6.Paragraphs(5,1).IndentLevel = 2
$tr6.Paragraphs(7,1).IndentLevel = 2
$tr6.Paragraphs(8,1).IndentLevel = 2
$tr6.Paragraphs(9,1).IndentLevel = 2

$para6_1 = $tr6.Paragraphs(1,1)
$para6_1.Characters(1,14).Text = "Time spent on "
$para6_1.Characters(15,6).Text = "design"

$para6_8 = $tr6.Paragraphs(8,1)
$para6_8.Characters(1,74).Text = "How to decode message in Message Handler Class for use in Driver class to "
$para6_8.Characters(75,17).Text = "implement action"
